# Test plan workbook - report + Collection dropdown work
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text values, entered in the order that reproduces the shared-string pool ---
$ws.Range("A2").Value = "Pend"

$ws.Range("C1").Value = "CIPC Status"
$ws.Range("D1").Value = "Account Status"
$ws.Range("E1").Value = "TT"
$ws.Range("F1").Value = "Number called"
$ws.Range("G1").Value = "Person contacted"
$ws.Range("H1").Value = "Email used"
$ws.Range("I1").Value = "PTP Date"
$ws.Range("J1").Value = "PTP Amount"
$ws.Range("K1").Value = "Outcome resolution"
$ws.Range("L1").Value = "Debit date"
$ws.Range("M1").Value = "Debit amount"
$ws.Range("N1").Value = "Pend reason"
$ws.Range("O1").Value = "Outcome notes"
$ws.Range("P1").Value = "NVDT"
$ws.Range("Q1").Value = "Assignment"
$ws.Range("R1").Value = "Next steps"

$ws.Range("C2").Value = "Final dereg"
$ws.Range("D2").Value = "Cancelled"
$ws.Range("E2").Value = "Admin"
$ws.Range("G2").Value = "MR XX"
$ws.Range("I2").Value = "10/30/2020 12:00 AM"
$ws.Range("K2").Value = "callback"
$ws.Range("O2").Value = "no notes because"
$ws.Range("Q2").Value = "Me"
$ws.Range("R2").Value = "Big ones"

$ws.Range("B1").Value = "Account Number"

$ws.Range("B3").Value = "BAC101"
$ws.Range("C3").Value = "Option 2"
$ws.Range("D3").Value = "Open"
$ws.Range("G3").Value = "Joe Soap"
$ws.Range("R3").Value = "nnnnn"
$ws.Range("P3").Value = "10/23/2020 2:00 AM"
$ws.Range("O3").Value = "dgdgdgdgdgdggdgdgdgd"
$ws.Range("K3").Value = "dispute"
$ws.Range("E3").Value = "Call"

# --- Numeric values ---
$ws.Range("J2").Value = 1234
$ws.Range("N2").Value = 1
$ws.Range("P2").Value = 43962.208333333336
$ws.Range("F3").Value = 1112345678
$ws.Range("L3").Value = 43962
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 3

# --- Formats: date/time number format on L3 (date only) first ---
$ws.Range("L3").NumberFormat = "m/d/yy h:mm"

# --- Fills: green highlight block (row 2 outcome block + D3:E3) ---
$green = 5296274   # RGB(146, 208, 80)
$red = 255         # RGB(255, 0, 0)

$ws.Range("C2:E2").Interior.Color = $green
$ws.Range("F2:H2").Interior.Color = $green
$ws.Range("I2:K2").Interior.Color = $green
$ws.Range("L2:N2").Interior.Color = $green
$ws.Range("O2").Interior.Color = $green
$ws.Range("R2").Interior.Color = $green
$ws.Range("D3:E3").Interior.Color = $green

# --- P2: green fill AND date format together (creates the combined style) ---
$ws.Range("P2").Interior.Color = $green
$ws.Range("P2").NumberFormat = "m/d/yy h:mm"

# --- C3: red highlight ---
$ws.Range("C3").Interior.Color = $red

# --- F2, H2, L2, M2 stay blank but keep the green fill applied above ---

# --- Column widths (best-fit, matching the column's widest entry) ---
$ws.Columns.Item(1).ColumnWidth = 4.330729166666667
$ws.Columns.Item(2).ColumnWidth = 13.998697916666666
$ws.Columns.Item(3).ColumnWidth = 9.498697916666666
$ws.Columns.Item(4).ColumnWidth = 12.498697916666666
$ws.Columns.Item(5).ColumnWidth = 5.666666666666667
$ws.Columns.Item(6).ColumnWidth = 12.330729166666666
$ws.Columns.Item(7).ColumnWidth = 14.330729166666666
$ws.Columns.Item(8).ColumnWidth = 9.330729166666666
$ws.Columns.Item(9).ColumnWidth = 18.498697916666668
$ws.Columns.Item(10).ColumnWidth = 10.330729166666666
$ws.Columns.Item(11).ColumnWidth = 16.666666666666668
$ws.Columns.Item(12).ColumnWidth = 11.998697916666666
$ws.Columns.Item(13).ColumnWidth = 11.498697916666666
$ws.Columns.Item(14).ColumnWidth = 10.330729166666666
$ws.Columns.Item(15).ColumnWidth = 14.498697916666666
$ws.Columns.Item(16).ColumnWidth = 17.498697916666668
$ws.Columns.Item(18).ColumnWidth = 8.998697916666666

# --- Selection ---
[void]$ws.Range("E3").Select()
